$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column B text values
$ws.Range("B2").Value = "<may>"
$ws.Range("B5").Value = "<ab>"
$ws.Range("B9").Value = "<brov>"

# Update column C numeric values
$ws.Range("C2").Value = 16
$ws.Range("C3").Value = 18
$ws.Range("C4").Value = 13
$ws.Range("C5").Value = 17
$ws.Range("C6").Value = 17
$ws.Range("C7").Value = 13
$ws.Range("C8").Value = 16
$ws.Range("C9").Value = 22
$ws.Range("C10").Value = 14
$ws.Range("C11").Value = 16
$ws.Range("C12").Value = 15
$ws.Range("C13").Value = 14
$ws.Range("C14").Value = 16
$ws.Range("C15").Value = 21
$ws.Range("C16").Value = 16
$ws.Range("C17").Value = 19
$ws.Range("C18").Value = 14
